$d = $word.ActiveDocument

$d.Content.Find.Execute("432×5=2160", $true, $false, $false, $false, $false, $true, 1, $false, "685×2=1370", 2)
$d.Content.Find.Execute("312×8=2496", $true, $false, $false, $false, $false, $true, 1, $false, "481×4=1924", 2)
$d.Content.Find.Execute("240×3=720", $true, $false, $false, $false, $false, $true, 1, $false, "951×8=7608", 2)
$d.Content.Find.Execute("374×2=748", $true, $false, $false, $false, $false, $true, 1, $false, "608×3=1824", 2)
$d.Content.Find.Execute("700×3=2100", $true, $false, $false, $false, $false, $true, 1, $false, "329×5=1645", 2)
$d.Content.Find.Execute("155×7=1085", $true, $false, $false, $false, $false, $true, 1, $false, "973×2=1946", 2)
$d.Content.Find.Execute("635×2=1270", $true, $false, $false, $false, $false, $true, 1, $false, "150×9=1350", 2)
$d.Content.Find.Execute("681×7=4767", $true, $false, $false, $false, $false, $true, 1, $false, "663×8=5304", 2)
$d.Content.Find.Execute("381×2=762", $true, $false, $false, $false, $false, $true, 1, $false, "827×8=6616", 2)
$d.Content.Find.Execute("273×8=2184", $true, $false, $false, $false, $false, $true, 1, $false, "255×3=765", 2)
$d.Content.Find.Execute("431×3=1293", $true, $false, $false, $false, $false, $true, 1, $false, "264×3=792", 2)
$d.Content.Find.Execute("333×9=2997", $true, $false, $false, $false, $false, $true, 1, $false, "641×6=3846", 2)
$d.Content.Find.Execute("958×2=1916", $true, $false, $false, $false, $false, $true, 1, $false, "891×3=2673", 2)
$d.Content.Find.Execute("565×4=2260", $true, $false, $false, $false, $false, $true, 1, $false, "169×7=1183", 2)
$d.Content.Find.Execute("835×7=5845", $true, $false, $false, $false, $false, $true, 1, $false, "348×4=1392", 2)
$d.Content.Find.Execute("704×6=4224", $true, $false, $false, $false, $false, $true, 1, $false, "377×5=1885", 2)
$d.Content.Find.Execute("229×5=1145", $true, $false, $false, $false, $false, $true, 1, $false, "840×5=4200", 2)
$d.Content.Find.Execute("931×3=2793", $true, $false, $false, $false, $false, $true, 1, $false, "956×9=8604", 2)
$d.Content.Find.Execute("751×3=2253", $true, $false, $false, $false, $false, $true, 1, $false, "994×2=1988", 2)
$d.Content.Find.Execute("564×4=2256", $true, $false, $false, $false, $false, $true, 1, $false, "712×7=4984", 2)
$d.Content.Find.Execute("109×2=218", $true, $false, $false, $false, $false, $true, 1, $false, "671×6=4026", 2)
$d.Content.Find.Execute("865×5=4325", $true, $false, $false, $false, $false, $true, 1, $false, "276×5=1380", 2)
$d.Content.Find.Execute("588×3=1764", $true, $false, $false, $false, $false, $true, 1, $false, "951×2=1902", 2)
$d.Content.Find.Execute("907×3=2721", $true, $false, $false, $false, $false, $true, 1, $false, "150×4=600", 2)
$d.Content.Find.Execute("788×4=3152", $true, $false, $false, $false, $false, $true, 1, $false, "473×2=946", 2)
